# v7.3: many updates according to changelog
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the gear ratio base value (A19); re-enter the B2:B16 formula as a
# single range assignment so Excel stores it as a shared formula group,
# and this cascades/recalculates the dependent ratio cells B2:B16.
$ws.Range("A19").Value = 1.5
$ws.Range("B2:B16").Formula = "=B1*A$19"

# Update the active selection to match the author's last cursor position
$ws.Range("E20").Select()

$wb.Application.Calculate()
